$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 1.069761271434882
$ws.Range("D2").Value2 = 1.067785388946211
$ws.Range("E2").Value2 = 1.073269561572123
$ws.Range("F2").Value2 = 1.082819363328977
$ws.Range("J2").Value2 = 1.07469435766929
$ws.Range("K2").Value2 = 1.070492320365496
$ws.Range("L2").Value2 = 1.075961869338507
$ws.Range("M2").Value2 = 1.085486588178515
$ws.Range("N2").Value2 = 1.076220546926347
$ws.Range("C3").Value2 = 1.073784900911161
$ws.Range("D3").Value2 = 1.071509064804204
$ws.Range("E3").Value2 = 1.076957498081294
$ws.Range("F3").Value2 = 1.08677103995097
$ws.Range("J3").Value2 = 1.078360888699483
$ws.Range("K3").Value2 = 1.074024634390174
$ws.Range("L3").Value2 = 1.079459663349572
$ws.Range("M3").Value2 = 1.089249426160168
$ws.Range("N3").Value2 = 1.079892284851159
$ws.Range("C4").Value2 = 1.076362090783612
$ws.Range("D4").Value2 = 1.07389328218832
$ws.Range("E4").Value2 = 1.079318663033244
$ws.Range("F4").Value2 = 1.089302230544027
$ws.Range("J4").Value2 = 1.080707745567634
$ws.Range("K4").Value2 = 1.076285050242309
$ws.Range("L4").Value2 = 1.081697799192504
$ws.Range("M4").Value2 = 1.091658474222222
$ws.Range("N4").Value2 = 1.082242474525251
$ws.Range("C5").Value2 = 1.07743944674907
$ws.Range("D5").Value2 = 1.074889764089809
$ws.Range("E5").Value2 = 1.080305469383535
$ws.Range("F5").Value2 = 1.09036037608612
$ws.Range("J5").Value2 = 1.081688428980049
$ws.Range("K5").Value2 = 1.077229484796222
$ws.Range("L5").Value2 = 1.082632879130794
$ws.Range("M5").Value2 = 1.092665275982179
$ws.Range("N5").Value2 = 1.083224550620597
$ws.Range("C6").Value2 = 1.077619988369249
$ws.Range("D6").Value2 = 1.075056740895711
$ws.Range("E6").Value2 = 1.080470822488937
$ws.Range("F6").Value2 = 1.090537699503248
$ws.Range("J6").Value2 = 1.081852747741762
$ws.Range("K6").Value2 = 1.0773877223674
$ws.Range("L6").Value2 = 1.082789546670649
$ws.Range("M6").Value2 = 1.092833978650146
$ws.Range("N6").Value2 = 1.0833891027338
$ws.Range("C7").Value2 = 1.076376510130171
$ws.Range("D7").Value2 = 1.073906619924005
$ws.Range("E7").Value2 = 1.079331871424508
$ws.Range("F7").Value2 = 1.089316392706307
$ws.Range("J7").Value2 = 1.080720872565408
$ws.Range("K7").Value2 = 1.076297692531983
$ws.Range("L7").Value2 = 1.081710316440782
$ws.Range("M7").Value2 = 1.091671950315746
$ws.Range("N7").Value2 = 1.082255620164867
$ws.Range("C8").Value2 = 1.071126668725626
$ws.Range("D8").Value2 = 1.069049174327533
$ws.Range("E8").Value2 = 1.074521252353222
$ws.Range("F8").Value2 = 1.084160324290596
$ws.Range("J8").Value2 = 1.075938910077688
$ws.Range("K8").Value2 = 1.071691425488777
$ws.Range("L8").Value2 = 1.077149295479844
$ws.Range("M8").Value2 = 1.086763716044525
$ws.Range("N8").Value2 = 1.077466866741916
$ws.Range("C9").Value2 = 1.061663718938187
$ws.Range("D9").Value2 = 1.060286977572874
$ws.Range("E9").Value2 = 1.065842217899969
$ws.Range("F9").Value2 = 1.074867182822569
$ws.Range("J9").Value2 = 1.067306791918153
$ws.Range("K9").Value2 = 1.063372311689238
$ws.Range("L9").Value2 = 1.068910423956395
$ws.Range("M9").Value2 = 1.077907974389677
$ws.Range("N9").Value2 = 1.068822489984478
$ws.Range("C10").Value2 = 1.05519785427093
$ws.Range("D10").Value2 = 1.054295611293698
$ws.Range("E10").Value2 = 1.059906854158195
$ws.Range("F10").Value2 = 1.068518011982954
$ws.Range("J10").Value2 = 1.061400184710705
$ws.Range("K10").Value2 = 1.057677090588838
$ws.Range("L10").Value2 = 1.063269155963709
$ws.Range("M10").Value2 = 1.07185131371572
$ws.Range("N10").Value2 = 1.062907494717297
$ws.Range("C11").Value2 = 1.052357348190431
$ws.Range("D11").Value2 = 1.051662558193424
$ws.Range("E11").Value2 = 1.057298213136341
$ws.Range("F11").Value2 = 1.065728960723681
$ws.Range("J11").Value2 = 1.058803370737038
$ws.Range("K11").Value2 = 1.055172552928916
$ws.Range("L11").Value2 = 1.060788115781958
$ws.Range("M11").Value2 = 1.069189243321587
$ws.Range("N11").Value2 = 1.060306992969929
$ws.Range("C12").Value2 = 1.051295830715829
$ws.Range("D12").Value2 = 1.050678419744675
$ws.Range("E12").Value2 = 1.05632316946924
$ws.Range("F12").Value2 = 1.064686703683364
$ws.Range("J12").Value2 = 1.057832623876866
$ws.Range("K12").Value2 = 1.054236202368447
$ws.Range("L12").Value2 = 1.059860515943885
$ws.Range("M12").Value2 = 1.068194210348589
$ws.Range("N12").Value2 = 1.05933486753787
$ws.Range("C13").Value2 = 1.051523825717757
$ws.Range("D13").Value2 = 1.050889801750513
$ws.Range("E13").Value2 = 1.056532599346114
$ws.Range("F13").Value2 = 1.064910560422093
$ws.Range("J13").Value2 = 1.058041136471595
$ws.Range("K13").Value2 = 1.054437331228016
$ws.Range("L13").Value2 = 1.060059766676769
$ws.Range("M13").Value2 = 1.068407934598043
$ws.Range("N13").Value2 = 1.059543676244399
$ws.Range("C14").Value2 = 1.052269735750746
$ws.Range("D14").Value2 = 1.05158133523092
$ws.Range("E14").Value2 = 1.057217741383966
$ws.Range("F14").Value2 = 1.06564293728982
$ws.Range("J14").Value2 = 1.058723256138749
$ws.Range("K14").Value2 = 1.055095279023089
$ws.Range("L14").Value2 = 1.06071156473617
$ws.Range("M14").Value2 = 1.069107122210939
$ws.Range("N14").Value2 = 1.060226764599721
$ws.Range("C15").Value2 = 1.052728454518375
$ws.Range("D15").Value2 = 1.052006594134751
$ws.Range("E15").Value2 = 1.057639065973708
$ws.Range("F15").Value2 = 1.066093337684767
$ws.Range("J15").Value2 = 1.05914270571262
$ws.Range("K15").Value2 = 1.055499851784162
$ws.Range("L15").Value2 = 1.061112351524428
$ws.Range("M15").Value2 = 1.069537081545365
$ws.Range("N15").Value2 = 1.060646809840099
$ws.Range("C16").Value2 = 1.055385485092562
$ws.Range("D16").Value2 = 1.054469518008764
$ws.Range("E16").Value2 = 1.060079144370747
$ws.Range("F16").Value2 = 1.06870224815731
$ws.Range("J16").Value2 = 1.061571676634513
$ws.Range("K16").Value2 = 1.057842474877156
$ws.Range("L16").Value2 = 1.063432983854576
$ws.Range("M16").Value2 = 1.072027130160851
$ws.Range("N16").Value2 = 1.063079230179307
$ws.Range("C17").Value2 = 1.057041050305958
$ws.Range("D17").Value2 = 1.056003874084444
$ws.Range("E17").Value2 = 1.061599215428059
$ws.Range("F17").Value2 = 1.070327882768646
$ws.Range("J17").Value2 = 1.063084611336228
$ws.Range("K17").Value2 = 1.059301451397464
$ws.Range("L17").Value2 = 1.064878204044377
$ws.Range("M17").Value2 = 1.073578298551096
$ws.Range("N17").Value2 = 1.064594313421834
$ws.Range("C18").Value2 = 1.058002803626233
$ws.Range("D18").Value2 = 1.056895118798757
$ws.Range("E18").Value2 = 1.062482143439307
$ws.Range("F18").Value2 = 1.071272266771258
$ws.Range("J18").Value2 = 1.063963316726029
$ws.Range("K18").Value2 = 1.060148755128849
$ws.Range("L18").Value2 = 1.065717496842224
$ws.Range("M18").Value2 = 1.074479278086787
$ws.Range("N18").Value2 = 1.065474266674078
$ws.Range("C19").Value2 = 1.058330082532018
$ws.Range("D19").Value2 = 1.057198387761434
$ws.Range("E19").Value2 = 1.062782579173821
$ws.Range("F19").Value2 = 1.07159363797398
$ws.Range("J19").Value2 = 1.064262302532386
$ws.Range("K19").Value2 = 1.060437045600964
$ws.Range("L19").Value2 = 1.066003057911242
$ws.Range("M19").Value2 = 1.074785854370917
$ws.Range("N19").Value2 = 1.065773677074574
$ws.Range("C20").Value2 = 1.056863830407798
$ws.Range("D20").Value2 = 1.055839638902445
$ws.Range("E20").Value2 = 1.061436511270897
$ws.Range("F20").Value2 = 1.070153864904629
$ws.Range("J20").Value2 = 1.062922679076663
$ws.Range("K20").Value2 = 1.059145300908941
$ws.Range("L20").Value2 = 1.064723528128731
$ws.Range("M20").Value2 = 1.073412266995192
$ws.Range("N20").Value2 = 1.064432151199888
$ws.Range("C21").Value2 = 1.052050263780349
$ws.Range("D21").Value2 = 1.051377866732292
$ws.Range("E21").Value2 = 1.057016154237538
$ws.Range("F21").Value2 = 1.065427446346162
$ws.Range("J21").Value2 = 1.058522561683257
$ws.Range("K21").Value2 = 1.054901699174301
$ws.Range("L21").Value2 = 1.060519795184712
$ws.Range("M21").Value2 = 1.068901403005153
$ws.Range("N21").Value2 = 1.060025785135082
$ws.Range("C22").Value2 = 1.048986436919784
$ws.Range("D22").Value2 = 1.048537100465998
$ws.Range("E22").Value2 = 1.054201585194992
$ws.Range("F22").Value2 = 1.062419272844666
$ws.Range("J22").Value2 = 1.055720159730902
$ws.Range("K22").Value2 = 1.052198408252214
$ws.Range("L22").Value2 = 1.057841704016116
$ws.Range("M22").Value2 = 1.066029096097795
$ws.Range("N22").Value2 = 1.057219403450515
$ws.Range("C23").Value2 = 1.050614275254111
$ws.Range("D23").Value2 = 1.05004650479813
$ws.Range("E23").Value2 = 1.055697086010732
$ws.Range("F23").Value2 = 1.06401752373338
$ws.Range("J23").Value2 = 1.057209264314162
$ws.Range("K23").Value2 = 1.053634902454558
$ws.Range("L23").Value2 = 1.059264825958131
$ws.Range("M23").Value2 = 1.067555286184471
$ws.Range("N23").Value2 = 1.058710622733086
$ws.Range("C24").Value2 = 1.056943920512944
$ws.Range("D24").Value2 = 1.055913861185715
$ws.Range("E24").Value2 = 1.061510041702205
$ws.Range("F24").Value2 = 1.070232507874907
$ws.Range("J24").Value2 = 1.062995860906221
$ws.Range("K24").Value2 = 1.059215869985502
$ws.Range("L24").Value2 = 1.064793430869574
$ws.Range("M24").Value2 = 1.073487301203202
$ws.Range("N24").Value2 = 1.064505436956037
$ws.Range("C25").Value2 = 1.064136719599562
$ws.Range("D25").Value2 = 1.062577606536118
$ws.Range("E25").Value2 = 1.068111258050262
$ws.Range("F25").Value2 = 1.077295700824986
$ws.Range("J25").Value2 = 1.069564138732758
$ws.Range("K25").Value2 = 1.065548292893455
$ws.Range("L25").Value2 = 1.071065588445654
$ws.Range("M25").Value2 = 1.080223287716545
$ws.Range("N25").Value2 = 1.07108304249048
